# Updates the result table on the active sheet (pl_mw.xlsx equivalent):
# columns B, C, E, F, G, K, O for data rows 2-25 get new computed values
# (per commit "case with 380 kV done").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.3711956348603849,0.06581034266109498,0.4155111673881606,0.4443680307746121,0.002425793169282033,0.3605865209266312,2.241390511598453),
    @(0.328673998069462,0.06023286219627266,0.3625780745220482,0.3878228170618172,0.002428378515667693,0.3151912533715517,2.269252506580187),
    @(0.3025654375947795,0.05678764258971114,0.3301608604842272,0.3531389305169483,0.002430047342677219,0.2872566965200463,2.288147823102989),
    @(0.2919263907652692,0.05537856450142442,0.3169700923546657,0.3390132514313251,0.002430747940476173,0.2758580777435213,2.296296152141238),
    @(0.2901598222964594,0.05514428061064791,0.3147809073601451,0.336668177824194,0.002430865516651234,0.2739644495740947,2.297676217467867),
    @(0.3024219534814847,0.05676865992923297,0.3299828882651923,0.3529483938344953,0.002430056707831238,0.2871030311368372,2.288255900709899),
    @(0.3565344390181906,0.06389154376331874,0.3972414439682552,0.4248636149813478,0.002426667738725298,0.3449472884047395,2.250625484100723),
    @(0.4626346507122605,0.07769377383846177,0.5298791984712778,0.5661985755041457,0.002420664961906481,0.457876284631368,2.191075715285237),
    @(0.540567927228949,0.08773145745531963,0.6278987771391087,0.6702781546542269,0.002416642476409051,0.5405279446148938,2.156086578636888),
    @(0.5760162539478131,0.0922752346481559,0.6726402600709349,0.7176906081379002,0.002414895839102381,0.5780578516655339,2.142088826596819),
    @(0.5894387516663926,0.09399257589083732,0.6896064642672712,0.7356546913071611,0.002414246330527422,0.5922592537166054,2.137065606750099),
    @(0.5865480288642573,0.09362286301222866,0.68595141229234,0.7317853510981394,0.002414385685302515,0.5892011958850105,2.138135084872687),
    @(0.5771205538084132,0.09241658762280736,0.6740355974109491,0.7191683204515869,0.002414842165428857,0.5792264198615271,2.141669994781864),
    @(0.5713458077255211,0.09167727893127164,0.6667399387604718,0.7114413442032514,0.00241512332117616,0.5731152185174437,2.143871402053065),
    @(0.5382511782444794,0.08743405426891115,0.6249780202225139,0.6671810134426437,0.002416758292597632,0.5380738529248958,2.157040105675151),
    @(0.5179473698090078,0.08482518544407469,0.599398566133587,0.6400460337125793,0.002417782563722382,0.5165591632076314,2.165611225891411),
    @(0.5062688146661856,0.08332252708177634,0.584700132128873,0.6244449056556647,0.002418379533240251,0.5041780397451134,2.170721642454851),
    @(0.5023146117616193,0.08281339307686153,0.5797258797249896,0.619163680173358,0.002418583004403862,0.4999849120255533,2.172482903818121),
    @(0.5201087835192766,0.08510312242636076,0.6021200614348174,0.642933953830422,0.002417672717626329,0.5188501089327815,2.164680120122767),
    @(0.5798896628564592,0.09277098971040232,0.6775349073391368,0.7228739723491628,0.002414707763596837,0.5821565382581753,2.140624165012042),
    @(0.6189538058017661,0.09776318648991378,0.7269611030361602,0.7751780083420101,0.002412839358399942,0.6234703529579519,2.126519943063101),
    @(0.5981052450352138,0.09510053496903481,0.700568185997156,0.7472568307830727,0.00241383023469105,0.601426091770179,2.1338991191638),
    @(0.5191316258707275,0.08497747578795156,0.6008896501042926,0.6416283278902171,0.002417722353809787,0.5178144098461814,2.165100503823794),
    @(0.4339342978255729,0.07397783566773342,0.4939054811941048,0.5279251897347166,0.002422220474494496,0.4273809917877713,2.205652514632888),
)

$cols = @(2, 3, 5, 6, 7, 11, 15)   # B, C, E, F, G, K, O

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowValues[$j]
    }
}
